$d = $word.ActiveDocument

# 1. Update the date paragraph from "10/05/23" to "2023-05-10"
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*10/05/23*") {
        $p.Range.Text = "2023-05-10"
        break
    }
}

# 2. Trim the TOC heading paragraph from "Indice del contenuto" to just "Indice"
foreach ($p in $d.Paragraphs) {
    $full = $p.Range.Text
    $target = " del contenuto"
    if ($full -like "*$target*") {
        $idx = $full.IndexOf($target)
        $start = $p.Range.Start
        $r = $d.Range($start + $idx, $start + $idx + $target.Length)
        $r.Delete()
        break
    }
}
